$d = $word.ActiveDocument

# 1. Apply justified ("both") alignment to paragraphs 1 through 84
#    (from the document title through the last bullet of "Propuesta
#    tecnológica 2" -- "Capacidades para el envio de mensajes...").
for ($i = 1; $i -le 84; $i++) {
    $d.Paragraphs($i).Range.ParagraphFormat.Alignment = 3
}

# 2. Split "Manejo de fabrica de archivos, imágenes, etc." so "fabrica"
#    is wrapped in spell-check proofErr markers.
$p82 = $d.Paragraphs(82)
$xml82 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="NoSpacing"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr><w:jc w:val="both"/></w:pPr><w:r><w:t xml:space="preserve">Manejo de </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>fabrica</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> de archivos, im&#225;genes, etc.</w:t></w:r></w:p>
'@
$p82.Range.InsertXML($xml82) | Out-Null

# 3. Split "Para la fase de diseño ... asi como ..." so "asi" is wrapped
#    in spell-check proofErr markers.
$p104 = $d.Paragraphs(104)
$xml104 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Para la fase de dise&#241;o se realizan los diagramas relacionales de base de datos </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>asi</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> como los diagramas de clases del aplicativo. </w:t></w:r></w:p>
'@
$p104.Range.InsertXML($xml104) | Out-Null

# 4. Split "Muck Ups de un posible producto terminado." so "Muck" is
#    wrapped in spell-check proofErr markers.
$p110 = $d.Paragraphs(110)
$xml110 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:lastRenderedPageBreak/><w:t>Muck</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> Ups de un posible producto terminado.</w:t></w:r></w:p>
'@
$p110.Range.InsertXML($xml110) | Out-Null
